# =====================================================================
# Commit: "add file ngu phap"
# Adds two new worksheets ("ngu phap" and "loai tu") with study-notes
# content, and updates the view-state (selection) on the two existing
# worksheets ("tong" and "nghe").
# =====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Add the two new worksheets at the end of the tab strip
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsGrammar = $wb.Worksheets.Add($null, $lastSheet)
$wsGrammar.Name = "ngu phap"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsWordType = $wb.Worksheets.Add($null, $lastSheet2)
$wsWordType.Name = "loai tu"

# ---------------------------------------------------------------
# 2) "ngu phap" sheet — step-by-step fill-in-the-blank strategy notes
# ---------------------------------------------------------------

$wsGrammar.Range("A2").Value = "CÁC BƯỚC LÀM BÀI DẠNG ĐIỀN TỪ TRONG ĐOẠN VĂN"
$wsGrammar.Range("A2:F2").Interior.Color = 65535

$wsGrammar.Range("A4").Value = "BƯỚC 1: Đọc câu trước – câu sau của chỗ trống"
$wsGrammar.Range("A4").Interior.Color = 65535

$wsGrammar.Range("B6").Value = "Không cần đọc toàn bài ngay."

$wsGrammar.Range("B7").Value = "Chỉ nhìn 1–2 câu xung quanh ô trống, để biết nghĩa chung."
$wsGrammar.Range("B7").Characters(10, 26).Font.Bold = $true

$wsGrammar.Range("A9").Value = "BƯỚC 2: Xác định loại từ cần điền"
$wsGrammar.Range("A9").Interior.Color = 65535

$wsGrammar.Range("B11").Value = "Dựa vào:"

$wsGrammar.Range("B12").HorizontalAlignment = -4131
$wsGrammar.Range("B12").VerticalAlignment = -4108
$wsGrammar.Range("B12").IndentLevel = 1

$wsGrammar.Range("B13").Value = "phía trước có trợ từ gì?"
$wsGrammar.Range("B13").HorizontalAlignment = -4131
$wsGrammar.Range("B13").VerticalAlignment = -4108
$wsGrammar.Range("B13").IndentLevel = 1

$wsGrammar.Range("B14").HorizontalAlignment = -4131
$wsGrammar.Range("B14").VerticalAlignment = -4108
$wsGrammar.Range("B14").IndentLevel = 1

$wsGrammar.Range("B15").Value = "phía sau là danh từ hay động từ?"
$wsGrammar.Range("B15").HorizontalAlignment = -4131
$wsGrammar.Range("B15").VerticalAlignment = -4108
$wsGrammar.Range("B15").IndentLevel = 1

$wsGrammar.Range("B16").HorizontalAlignment = -4131
$wsGrammar.Range("B16").VerticalAlignment = -4108
$wsGrammar.Range("B16").IndentLevel = 1

$wsGrammar.Range("B17").Value = "có “とき”, “ので”, “から”… không?"
$wsGrammar.Range("B17").HorizontalAlignment = -4131
$wsGrammar.Range("B17").VerticalAlignment = -4108
$wsGrammar.Range("B17").IndentLevel = 1

$wsGrammar.Range("A19").Value = "BƯỚC 3: Loại trừ đáp án sai"
$wsGrammar.Range("A19").Interior.Color = 65535

$wsGrammar.Range("A21").Value = "BƯỚC 4: Kiểm tra nghĩa cho chắc chắn"
$wsGrammar.Range("A21").Interior.Color = 65535

$wsGrammar.Range("A23").Value = "BƯỚC 5: Đọc lại toàn đoạn để xem có bị “sai nghĩa” không"
$wsGrammar.Range("A23").Interior.Color = 65535

$wsGrammar.Range("G15").Select() | Out-Null

# ---------------------------------------------------------------
# 3) "loai tu" sheet — word-type cheat sheet table
# ---------------------------------------------------------------

$wsWordType.Columns.Item(1).ColumnWidth = 23.1
$wsWordType.Columns.Item(2).ColumnWidth = 37.1

$wsWordType.Range("A2").Value = "Loại từ"
$wsWordType.Range("B2").Value = "Dấu hiệu nhận biết"
$hdr = $wsWordType.Range("A2:B2")
$hdr.Font.Bold = $true
$hdr.Interior.Color = 65535
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4108
$hdr.WrapText = $true

$wsWordType.Range("A3").Value = "Động từ"
$wsWordType.Range("B3").Value = "có dạng ます, chia được thì, kết thúc bằng う/る"

$wsWordType.Range("A4").Value = "Danh từ"
$wsWordType.Range("B4").Value = "đứng trước trợ từ, sau の, không chia được"

$wsWordType.Range("A5").Value = "Tính từ -い"
$wsWordType.Range("B5").Value = "kết thúc bằng い, chia được"

$wsWordType.Range("A6").Value = "Tính từ -な"
$wsWordType.Range("B6").Value = "đứng trước N + な, không chia được"

$wsWordType.Range("A7").Value = "Trạng từ"
$wsWordType.Range("B7").Value = "kết thúc bằng に hoặc mô tả cách làm, đứng trước động từ"

$dataRowsA = $wsWordType.Range("A3:A7")
$dataRowsA.Font.Bold = $true
$dataRowsA.VerticalAlignment = -4108
$dataRowsA.WrapText = $true

$dataRowsB = $wsWordType.Range("B3:B7")
$dataRowsB.VerticalAlignment = -4108
$dataRowsB.WrapText = $true

$wsWordType.Rows.Item(3).RowHeight = 47.25
$wsWordType.Rows.Item(4).RowHeight = 47.25
$wsWordType.Rows.Item(5).RowHeight = 47.25
$wsWordType.Rows.Item(6).RowHeight = 47.25
$wsWordType.Rows.Item(7).RowHeight = 47.25

$wsWordType.Range("A9").Value = "Dùng JLPT Sensei (đề chuẩn nhất) để kiểm tra format."
$wsWordType.Range("A9").Font.Bold = $true
$wsWordType.Range("A9").Font.Size = 13.5
$wsWordType.Range("A9").VerticalAlignment = -4108
$wsWordType.Rows.Item(9).RowHeight = 18

$wsWordType.Range("A11").Value = "✔ Dùng Nihongo Pro để kiểm tra điểm thật của bạn."
$wsWordType.Range("A11").Font.Bold = $true
$wsWordType.Range("A11").Font.Size = 13.5
$wsWordType.Range("A11").VerticalAlignment = -4108
$wsWordType.Rows.Item(11).RowHeight = 18

$wsWordType.Range("F5").Select() | Out-Null
$wsWordType.Activate()

# ---------------------------------------------------------------
# 4) Update view-state on the pre-existing worksheets
# ---------------------------------------------------------------

$wsTong = $wb.Worksheets.Item("tong")
$wsTong.Range("B76").Select() | Out-Null

$wsNghe = $wb.Worksheets.Item("nghe")
$wsNghe.Range("F43").Select() | Out-Null

# "loai tu" stays the active / selected tab, matching activeTab="3"
$wsWordType.Activate()
